# colors.xlsx: add an "orange" entry, rename "blue" -> "blue -1",
# and append a new "Sheet2" with a "new sheet " note.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Write A3 before A2 so the shared-string table ends up in the same
# order as the target workbook: Red, "orange ", "blue -1", "new sheet ".
$ws1.Range("A3").Value = "orange "
$ws1.Range("A2").Value = "blue -1"

# Keep the same selection Excel leaves behind in the target file.
$ws1.Range("A3").Select()

# Add a new worksheet right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "new sheet "

# Leave Sheet1 as the active/selected sheet, matching the target file.
$ws1.Select()
